$d = $word.ActiveDocument

# 1. Remove the "Meta description" paragraph that follows the title heading.
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

# 2. Replace the final paragraph (the italic image-generation prompt) with two
#    paragraphs: a new bold "Play Eye of Horus Free: Review & Features" title,
#    followed by the (now relocated) meta-description text in italics.
$pLast = $d.Paragraphs.Item($d.Paragraphs.Count)
$r = $pLast.Range
$r.Collapse(1)

$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Eye of Horus Free: Review &amp; Features</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Experience the ancient Egyptian world of Eye of Horus slot game. Explore gameplay mechanics, features &amp; how to play it for free in our unbiased review.</w:t></w:r></w:p>
'@

$r.InsertXML($xml) | Out-Null
